$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.834.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5164"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06113"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.753.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07034"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6497"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.821.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006596"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.976.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.132"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.660"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.157"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.517"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.816"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08305"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.682"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.433"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04499"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9882"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6167"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.32%  "

$ws.Range("E38").Value = "  +2.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01586"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.936"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3860"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.970"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05410"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.301"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1127"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.658"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
